# Auto-generated edit script: applies numeric corrections to the
# per-Leve market-price / profit columns (H-N) across all 8 job sheets,
# per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1032.909
$ws.Range("I39").Value = 22.875
$ws.Range("K39").Value = 68.625
$ws.Range("M39").Value = 227.375
$ws.Range("H62").Value = 4955.8
$ws.Range("I62").Value = 4789.4
$ws.Range("K62").Value = 4789.4
$ws.Range("M62").Value = -4165.4
$ws.Range("H65").Value = 4955.8
$ws.Range("I65").Value = 4789.4
$ws.Range("K65").Value = 23947
$ws.Range("M65").Value = -20827
$ws.Range("H69").Value = 16583
$ws.Range("J69").Value = 16583
$ws.Range("L69").Value = 49749
$ws.Range("N69").Value = -51497
$ws.Range("H72").Value = 16583
$ws.Range("J72").Value = 16583
$ws.Range("L72").Value = 149247
$ws.Range("N72").Value = -157983
$ws.Range("H86").Value = 2467
$ws.Range("J86").Value = 1720.75
$ws.Range("L86").Value = 1720.75
$ws.Range("N86").Value = -3966.75
$ws.Range("H89").Value = 2467
$ws.Range("J89").Value = 1720.75
$ws.Range("L89").Value = 8603.75
$ws.Range("N89").Value = -19835.75
$ws.Range("H98").Value = 1577.1052
$ws.Range("I98").Value = 1577.1052
$ws.Range("K98").Value = 1577.1052
$ws.Range("M98").Value = -79.10519999999997
$ws.Range("H105").Value = 32999
$ws.Range("J105").Value = 32999
$ws.Range("L105").Value = 32999
$ws.Range("N105").Value = -39987
$ws.Range("H122").Value = 1577.1052
$ws.Range("I122").Value = 1577.1052
$ws.Range("K122").Value = 4731.3156
$ws.Range("M122").Value = -2281.3156
$ws.Range("H125").Value = 17746098
$ws.Range("I125").Value = 8474726
$ws.Range("J125").Value = 20836554
$ws.Range("K125").Value = 76272534
$ws.Range("L125").Value = 187528986
$ws.Range("M125").Value = -76270074
$ws.Range("N125").Value = -187533906
$ws.Range("H132").Value = 1870.125
$ws.Range("I132").Value = 1821.7778
$ws.Range("J132").Value = 2131.2
$ws.Range("K132").Value = 5465.3334
$ws.Range("L132").Value = 6393.599999999999
$ws.Range("M132").Value = -2935.3334
$ws.Range("N132").Value = -11453.6
$ws.Range("H135").Value = 55556308
$ws.Range("J135").Value = 333333540
$ws.Range("L135").Value = 3000001860
$ws.Range("N135").Value = -3000006930
$ws.Range("H138").Value = 3486.9
$ws.Range("I138").Value = 3127.1667
$ws.Range("J138").Value = 3600.5
$ws.Range("K138").Value = 9381.500100000001
$ws.Range("L138").Value = 10801.5
$ws.Range("M138").Value = -4241.500100000001
$ws.Range("N138").Value = -21081.5
$ws.Range("H140").Value = 99995
$ws.Range("J140").Value = 99995
$ws.Range("L140").Value = 99995
$ws.Range("N140").Value = -110355

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5228.03
$ws.Range("I32").Value = 2006.6274
$ws.Range("J32").Value = 15496.25
$ws.Range("K32").Value = 2006.6274
$ws.Range("L32").Value = 15496.25
$ws.Range("M32").Value = -1719.6274
$ws.Range("N32").Value = -16070.25
$ws.Range("H45").Value = 1365.7142
$ws.Range("I45").Value = 996.3333
$ws.Range("J45").Value = 2030.6
$ws.Range("K45").Value = 996.3333
$ws.Range("L45").Value = 2030.6
$ws.Range("M45").Value = -619.3333
$ws.Range("N45").Value = -2784.6
$ws.Range("H61").Value = 5070.5
$ws.Range("I61").Value = 4728.9
$ws.Range("K61").Value = 4728.9
$ws.Range("M61").Value = -4516.9
$ws.Range("H123").Value = 67101
$ws.Range("J123").Value = 67101
$ws.Range("L123").Value = 67101
$ws.Range("N123").Value = -76901
$ws.Range("H136").Value = 5070.5
$ws.Range("I136").Value = 4728.9
$ws.Range("K136").Value = 14186.7
$ws.Range("M136").Value = -11636.7

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 23809748
$ws.Range("I22").Value = 450
$ws.Range("K22").Value = 450
$ws.Range("M22").Value = -277
$ws.Range("H112").Value = 17110
$ws.Range("I112").Value = 17110
$ws.Range("K112").Value = 17110
$ws.Range("M112").Value = -15633

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 23558.7
$ws.Range("I58").Value = 26323.375
$ws.Range("K58").Value = 26323.375
$ws.Range("M58").Value = -26120.375
$ws.Range("H136").Value = 23558.7
$ws.Range("I136").Value = 26323.375
$ws.Range("K136").Value = 78970.125
$ws.Range("M136").Value = -76420.125

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 235538.4
$ws.Range("I81").Value = 294073.12
$ws.Range("J81").Value = 1399.5
$ws.Range("K81").Value = 882219.36
$ws.Range("L81").Value = 4198.5
$ws.Range("M81").Value = -881096.36
$ws.Range("N81").Value = -6444.5
$ws.Range("H84").Value = 235538.4
$ws.Range("I84").Value = 294073.12
$ws.Range("J84").Value = 1399.5
$ws.Range("K84").Value = 2646658.08
$ws.Range("L84").Value = 12595.5
$ws.Range("M84").Value = -2641042.08
$ws.Range("N84").Value = -23827.5
$ws.Range("H136").Value = 899
$ws.Range("I136").Value = 899
$ws.Range("K136").Value = 2697
$ws.Range("M136").Value = 2403

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 69997.5
$ws.Range("J111").Value = 69997.5
$ws.Range("L111").Value = 69997.5
$ws.Range("N111").Value = -76131.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 42499.5
$ws.Range("I45").Value = 35000
$ws.Range("J45").Value = 49999
$ws.Range("K45").Value = 35000
$ws.Range("L45").Value = 49999
$ws.Range("M45").Value = -34593
$ws.Range("N45").Value = -50813
$ws.Range("H104").Value = 55570.715
$ws.Range("J104").Value = 55570.715
$ws.Range("L104").Value = 55570.715
$ws.Range("N104").Value = -62558.715
$ws.Range("H132").Value = 25570044
$ws.Range("I132").Value = 32874628
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 98623884
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -98621354
$ws.Range("N132").Value = -17058.5
$ws.Range("H136").Value = 2998.5
$ws.Range("I136").Value = 2996.5
$ws.Range("J136").Value = 2999
$ws.Range("K136").Value = 8989.5
$ws.Range("L136").Value = 8997
$ws.Range("M136").Value = -6439.5
$ws.Range("N136").Value = -14097

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 46110.09
$ws.Range("J94").Value = 41721.2
$ws.Range("L94").Value = 41721.2
$ws.Range("N94").Value = -43523.2
$ws.Range("H100").Value = 1310.381
$ws.Range("I100").Value = 1374.7059
$ws.Range("K100").Value = 2749.4118
$ws.Range("M100").Value = -2208.4118
$ws.Range("H110").Value = 258999
$ws.Range("J110").Value = 258999
$ws.Range("L110").Value = 258999
$ws.Range("N110").Value = -267179
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H112").Value = 49346
$ws.Range("J112").Value = 49346
$ws.Range("L112").Value = 49346
$ws.Range("N112").Value = -52300
